# Natmi following Dr Hou advice
# Update computed LR-pair statistics for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value  = 3
$ws.Cells.Item(2, 7).Value  = 0.5587383333333333
$ws.Cells.Item(2, 8).Value  = 1.676215
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 138.2190853333333
$ws.Cells.Item(2, 14).Value = 414.657256
$ws.Cells.Item(2, 15).Value = 0.2249223651785973
$ws.Cells.Item(2, 16).Value = 0.2476599003709697
$ws.Cells.Item(2, 17).Value = 77.22830137400445
$ws.Cells.Item(2, 18).Value = 695.05471236604
$ws.Cells.Item(2, 19).Value = 0.2249223651785973
$ws.Cells.Item(2, 20).Value = 0.2476599003709697

# Row 3
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 7).Value  = 0.5587383333333333
$ws.Cells.Item(3, 8).Value  = 1.676215
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 147.91433
$ws.Cells.Item(3, 14).Value = 443.74299
$ws.Cells.Item(3, 15).Value = 0.2406993279341593
$ws.Cells.Item(3, 16).Value = 0.2650317656414439
$ws.Cells.Item(3, 17).Value = 82.64540622031666
$ws.Cells.Item(3, 18).Value = 743.8086559828499
$ws.Cells.Item(3, 19).Value = 0.2406993279341593
$ws.Cells.Item(3, 20).Value = 0.2650317656414439

# Row 4
$ws.Cells.Item(4, 5).Value  = 3
$ws.Cells.Item(4, 7).Value  = 0.5587383333333333
$ws.Cells.Item(4, 8).Value  = 1.676215
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 74.27261733333333
$ws.Cells.Item(4, 14).Value = 222.817852
$ws.Cells.Item(4, 15).Value = 0.1208629960061633
$ws.Cells.Item(4, 16).Value = 0.1330811078998542
$ws.Cells.Item(4, 17).Value = 41.49895842113111
$ws.Cells.Item(4, 18).Value = 373.49062579018
$ws.Cells.Item(4, 19).Value = 0.1208629960061633
$ws.Cells.Item(4, 20).Value = 0.1330811078998542

# Row 5
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 7).Value  = 0.5587383333333333
$ws.Cells.Item(5, 8).Value  = 1.676215
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 84.85695366666668
$ws.Cells.Item(5, 14).Value = 254.570861
$ws.Cells.Item(5, 15).Value = 0.138086767645209
$ws.Cells.Item(5, 16).Value = 0.1520460408212704
$ws.Cells.Item(5, 17).Value = 47.41283286345723
$ws.Cells.Item(5, 18).Value = 426.7154957711151
$ws.Cells.Item(5, 19).Value = 0.138086767645209
$ws.Cells.Item(5, 20).Value = 0.1520460408212704

# Row 6
$ws.Cells.Item(6, 5).Value  = 3
$ws.Cells.Item(6, 7).Value  = 0.5587383333333333
$ws.Cells.Item(6, 8).Value  = 1.676215
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 169.2560955
$ws.Cells.Item(6, 14).Value = 338.512191
$ws.Cells.Item(6, 15).Value = 0.275428543235871
$ws.Cells.Item(6, 16).Value = 0.2021811852664618
$ws.Cells.Item(6, 17).Value = 94.5698687061775
$ws.Cells.Item(6, 18).Value = 567.419212237065
$ws.Cells.Item(6, 19).Value = 0.275428543235871
$ws.Cells.Item(6, 20).Value = 0.2021811852664618
